# Apply the two logical changes described by the diff:
#  1. Header row height in both "herb_l_bin_2021" ranova tables grows
#     from 571 -> 637 twips (28.55pt -> 31.85pt), reflecting the rerun
#     with 1000 iterations.
#  2. The chi-square column header glyph "chi" got mangled into the
#     mojibake string "Ã\x8f\xE2\x80\xA1" (U+00CF U+2021) in both tables.

$d = $word.ActiveDocument

# --- 1. Fix header row heights (571 -> 637 twips == 28.55pt -> 31.85pt) ---
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    $row = $t.Rows.Item(1)
    if ([math]::Round($row.Height * 20) -eq 571) {
        $row.Height = 637 / 20
    }
}

# --- 2. Replace the chi character with the mojibake replacement text ---
$chi = [string]([char]0x03C7)
$mojibake = [string]([char]0x00CF) + [string]([char]0x2021)

$d.Content.Find.Execute($chi, $true, $false, $false, $false, $false, $true, 1, $false, $mojibake, 2) | Out-Null
